$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "kontenery produkcyjne firmy"
$ws.Range("D3").Value = "czesc stolarniana firmy"
$ws.Range("D4").Value = "plac poza stolarnia"
$ws.Range("D12").Value = "kontenery produkcyjne firmy"
$ws.Range("D13").Value = "hala cnc duzy maszyna "

$ws.Range("J6").Select()
